$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update question/answer text for the "is there a dining hall in the building" row
$ws.Range("A8").Value = "где в корпусе столовая"
$ws.Range("B8").Value = "Ахахах, о еде. `$`$удивление `$`$улыбка  В корпусе на первом этаже есть столовая и автомат с кофе, `$`$улыбка также прямо на остановке факультета есть киоск, выше по склону слата, и в библиотеке, где ты можешь ещё и почитать, есть столовая, где можно полноценно покушать."

# Update the view state: scroll back to top and move selection to B9
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("B9").Select()
